$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-32 with new values
$ws.Range("A2").Value = 46044
$ws.Range("B2").Value = 12.452
$ws.Range("C2").Value = 0
$ws.Range("A3").Value = 46044.01041666666
$ws.Range("B3").Value = 31.705
$ws.Range("C3").Value = 0
$ws.Range("A4").Value = 46044.02083333334
$ws.Range("B4").Value = 33.487
$ws.Range("C4").Value = 0
$ws.Range("A5").Value = 46044.03125
$ws.Range("B5").Value = 53.634
$ws.Range("C5").Value = 0
$ws.Range("A6").Value = 46044.04166666666
$ws.Range("B6").Value = 32.846
$ws.Range("C6").Value = 0
$ws.Range("A7").Value = 46044.05208333334
$ws.Range("B7").Value = 42.672
$ws.Range("C7").Value = 0
$ws.Range("A8").Value = 46044.0625
$ws.Range("B8").Value = 22.094
$ws.Range("C8").Value = 0
$ws.Range("A9").Value = 46044.07291666666
$ws.Range("B9").Value = 17.777
$ws.Range("C9").Value = 0
$ws.Range("A10").Value = 46044.08333333334
$ws.Range("B10").Value = 8.262
$ws.Range("C10").Value = 0
$ws.Range("A11").Value = 46044.09375
$ws.Range("B11").Value = 0.057
$ws.Range("C11").Value = 13.826
$ws.Range("A12").Value = 46044.10416666666
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 13.994
$ws.Range("A13").Value = 46044.11458333334
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 3.028
$ws.Range("A14").Value = 46044.125
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 10.455
$ws.Range("A15").Value = 46044.13541666666
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 3.913
$ws.Range("A16").Value = 46044.14583333334
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 2.943
$ws.Range("A17").Value = 46044.15625
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 3.704
$ws.Range("A18").Value = 46044.16666666666
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 6.617
$ws.Range("A19").Value = 46044.17708333334
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 25.675
$ws.Range("A20").Value = 46044.1875
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 5.335
$ws.Range("A21").Value = 46044.19791666666
$ws.Range("B21").Value = 1.067
$ws.Range("C21").Value = 3.25
$ws.Range("A22").Value = 46044.20833333334
$ws.Range("B22").Value = 7.292
$ws.Range("C22").Value = 0.138
$ws.Range("A23").Value = 46044.21875
$ws.Range("B23").Value = 32.407
$ws.Range("C23").Value = 0
$ws.Range("A24").Value = 46044.22916666666
$ws.Range("B24").Value = 12.031
$ws.Range("C24").Value = 0
$ws.Range("A25").Value = 46044.23958333334
$ws.Range("B25").Value = 1.364
$ws.Range("C25").Value = 3.929
$ws.Range("A26").Value = 46044.25
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 40.812
$ws.Range("A27").Value = 46044.26041666666
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 29.132
$ws.Range("A28").Value = 46044.27083333334
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 24.481
$ws.Range("A29").Value = 46044.28125
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 20.305
$ws.Range("A30").Value = 46044.29166666666
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 35.283
$ws.Range("A31").Value = 46044.30208333334
$ws.Range("B31").Value = 0.059
$ws.Range("C31").Value = 20.718
$ws.Range("A32").Value = 46044.32291666666
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0

# Remove now-unused rows 33-43 (shrinks the sheet dimension to A1:C32)
$ws.Range("A33:C43").EntireRow.Delete()
